$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - matches style of existing header row (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Data cells F2:F17 - time_taken values as text strings
$timeTaken = @(
    "2021-10-05 10:50:26.201919",
    "2021-10-05 10:50:26.201932",
    "2021-10-05 10:50:26.201937",
    "2021-10-05 10:50:26.201940",
    "2021-10-05 10:50:26.201944",
    "2021-10-05 10:50:26.201947",
    "2021-10-05 10:50:26.201950",
    "2021-10-05 10:50:26.201953",
    "2021-10-05 10:50:26.201957",
    "2021-10-05 10:50:26.201960",
    "2021-10-05 10:50:26.201963",
    "2021-10-05 10:50:26.201967",
    "2021-10-05 10:50:26.201970",
    "2021-10-05 10:50:26.201973",
    "2021-10-05 10:50:26.201976",
    "2021-10-05 10:50:26.201980"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timeTaken[$i]
}
